$wb = $excel.ActiveWorkbook

foreach ($name in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($name)
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value()
        if ($v -eq "Ready for handoff") {
            $cell.Value = "Handed back: in sync with en-US"
        }
    }
}

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("H2").Value = "2016-03-12 18:50:22"
$wsZh.Range("H3").Value = "2016-03-12 18:50:22"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("H2").Value = "2016-03-12 18:50:28"
$wsDe.Range("H3").Value = "2016-03-12 18:50:28"

# zh-cn Latest Target File (F) / Latest Handback File (G)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/6999aef4a2cfebfef8412ba27bb42cd5302f70b9/e2e/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.md", "", "", "92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2258437a52de4fa8b7da741ad685c297c204e0f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.1a7532fc668ee458d29b28e4d0235919d447ef1d.zh-cn.xlf", "", "", "92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.1a7532fc668ee458d29b28e4d0235919d447ef1d.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/6999aef4a2cfebfef8412ba27bb42cd5302f70b9/e2e/d0dc96df-64c6-47fb-94d3-fbc50197c361.md", "", "", "d0dc96df-64c6-47fb-94d3-fbc50197c361.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2258437a52de4fa8b7da741ad685c297c204e0f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d0dc96df-64c6-47fb-94d3-fbc50197c361.22ad50810204c9096ca26f259094112d8ab2d6ff.zh-cn.xlf", "", "", "d0dc96df-64c6-47fb-94d3-fbc50197c361.22ad50810204c9096ca26f259094112d8ab2d6ff.zh-cn.xlf")

# de-de Latest Target File (F) / Latest Handback File (G)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/6999aef4a2cfebfef8412ba27bb42cd5302f70b9/e2e/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.md", "", "", "92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc84705b8a1ef2a402362b46e0fbd6e9edeec32d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.1a7532fc668ee458d29b28e4d0235919d447ef1d.de-de.xlf", "", "", "92b3c6e3-c5d0-4ea3-8523-f7c514b343e6.1a7532fc668ee458d29b28e4d0235919d447ef1d.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/6999aef4a2cfebfef8412ba27bb42cd5302f70b9/e2e/d0dc96df-64c6-47fb-94d3-fbc50197c361.md", "", "", "d0dc96df-64c6-47fb-94d3-fbc50197c361.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc84705b8a1ef2a402362b46e0fbd6e9edeec32d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d0dc96df-64c6-47fb-94d3-fbc50197c361.22ad50810204c9096ca26f259094112d8ab2d6ff.de-de.xlf", "", "", "d0dc96df-64c6-47fb-94d3-fbc50197c361.22ad50810204c9096ca26f259094112d8ab2d6ff.de-de.xlf")
